$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.97"
$ws.Range("D2").Style = "Normal"

$ws.Range("D4").Value = "'5.339"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'0.05688"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'3.405"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'6.330"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'0.8127"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.9153"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'0.1409"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'0.07422"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'0.03106"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'0.03014"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'0.09362"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'3.720"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'0.001580"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "'0.04759"
$ws.Range("D17").Style = "Normal"

$ws.Range("E19").Value = "18OneONE"

$ws.Range("D20").Value = "'0.006438"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'0.004997"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").Value = "'0.001023"
$ws.Range("D22").Style = "Normal"

$ws.Range("D24").Value = "'3.698"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'2.146"
$ws.Range("D25").Style = "Normal"

$ws.Range("D40").Value = "'0.03980"
$ws.Range("D40").Style = "Normal"

$ws.Range("D43").Value = "'0.002710"
$ws.Range("D43").Style = "Normal"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

$ws.Range("D48").Value = "'0.2339"
$ws.Range("D48").Style = "Normal"
